$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 463. This shifts the existing rows 463:508
# down to 464:509 (matching the diff, which shows every row from 464..509
# now holding what used to be in the row above it).
$ws.Rows("463:463").Insert()

# Populate the newly inserted row 463 with the new weekly price record
# (same market/product/variety/quality as the row that used to be here,
# but a new date and updated volume/price figures).
$ws.Cells.Item(463, 1).Value2  = 4
$ws.Cells.Item(463, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(463, 3).Value2  = "Los Lagos"
$ws.Cells.Item(463, 4).Value2  = 44946
$ws.Cells.Item(463, 5).Value2  = 10
$ws.Cells.Item(463, 6).Value2  = "Fruta"
$ws.Cells.Item(463, 7).Value2  = 100102
$ws.Cells.Item(463, 8).Value2  = "Cítricos"
$ws.Cells.Item(463, 9).Value2  = 100102006
$ws.Cells.Item(463, 10).Value2 = "Pomelo"
$ws.Cells.Item(463, 11).Value2 = "Start Ruby"
$ws.Cells.Item(463, 12).Value2 = "Primera"
$ws.Cells.Item(463, 13).Value2 = 200
$ws.Cells.Item(463, 14).Value2 = 12000
$ws.Cells.Item(463, 15).Value2 = 13000
$ws.Cells.Item(463, 16).Value2 = 12500
$ws.Cells.Item(463, 17).Value2 = "$/caja 14 kilos empedrada"
$ws.Cells.Item(463, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(463, 19).Value2 = 893
$ws.Cells.Item(463, 20).Value2 = 14
